$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 4 (Patriotas vs Santa Fe) updates
$ws.Range("G4").Value = 3.3
$ws.Range("I4").Value = 2.3
$ws.Range("L4").Value = 3.2
$ws.Range("M4").Value = 1.13
$ws.Range("N4").Value = 6
$ws.Range("Z4").Value = 41
$ws.Range("AI4").Value = 10
$ws.Range("AJ4").Value = 21
$ws.Range("AP4").Value = 41
$ws.Range("AQ4").Value = 81
$ws.Range("AW4").Value = 4
$ws.Range("BA4").Value = 81

# Row 6 (America De Cali vs Ind. Medellin) updates
$ws.Range("H6").Value = 3.3
$ws.Range("J6").Value = 2.5
$ws.Range("L6").Value = 5.5
$ws.Range("Q6").Value = 2.5
$ws.Range("R6").Value = 1.5
$ws.Range("S6").Value = 1.53
$ws.Range("T6").Value = 2.38
$ws.Range("U6").Value = 2.25
$ws.Range("V6").Value = 1.57
$ws.Range("AC6").Value = 6.5
$ws.Range("AE6").Value = 21
$ws.Range("AF6").Value = 81
$ws.Range("AG6").Value = 10
$ws.Range("AI6").Value = 17
$ws.Range("AK6").Value = 41
$ws.Range("AN6").Value = 3.6
$ws.Range("AP6").Value = 26
$ws.Range("AQ6").Value = 34
$ws.Range("AS6").Value = 251
$ws.Range("AT6").Value = 2.38
$ws.Range("AU6").Value = 9.5
$ws.Range("AX6").Value = 29
$ws.Range("BA6").Value = 151

$wb.Save()
